# From v1.0.1 to v1.0.2
# The second test-step text (Description/Expected Results) of TC2 and TC4
# are swapped; TC3 stays where it is.
#
# Before:
#   TC2 step2: B20 = "Chefe Clica para realizar a autorização de pagamento."
#              D20 = "SYSTEM Apresenta a tela de Registrar Autorizações de Pagamento"
#   TC4 step2: B36 = "Chefe Dado um registro selecionado (...); e Clica para atribuir/desatribuir o registro a si mesmo."
#              D36 = "SYSTEM Atualiza a lista de registros de solicitações (...)."
#
# After:
#   TC2 step2: B20 = (old TC4 text)
#              D20 = (old TC4 result)
#   TC4 step2: B36 = (old TC2 text)
#              D36 = (old TC2 result)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc2StepText   = $ws.Range("B20").Value2
$tc2ResultText = $ws.Range("D20").Value2
$tc4StepText   = $ws.Range("B36").Value2
$tc4ResultText = $ws.Range("D36").Value2

$ws.Range("B20").Value = $tc4StepText
$ws.Range("D20").Value = $tc4ResultText
$ws.Range("B36").Value = $tc2StepText
$ws.Range("D36").Value = $tc2ResultText
